$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume table values (coinranking.com refresh)

$ws.Range("D2").Value = '42.699.91'
$ws.Range("E2").Value = '  -1.01%  '

$ws.Range("D3").Value = '2.370.11'
$ws.Range("E3").Value = '  +0.73%  '

$ws.Range("E4").Value = '  -0.13%  '

$ws.Range("D5").Value = '332.33'
$ws.Range("E5").Value = '  +7.70%  '

$ws.Range("D6").Value = '100.75'
$ws.Range("E6").Value = '  -6.38%  '

$ws.Range("D7").Value = '0.638'
$ws.Range("E7").Value = '  -0.08%  '

$ws.Range("E8").Value = '  +0.01%  '

$ws.Range("D9").Value = '0.637'
$ws.Range("E9").Value = '  +1.40%  '

$ws.Range("D10").Value = '40.07'
$ws.Range("E10").Value = '  -5.56%  '

$ws.Range("D11").Value = '0.0924'
$ws.Range("E11").Value = '  -1.40%  '

$ws.Range("D12").Value = '8.48'
$ws.Range("E12").Value = '  -4.52%  '

$ws.Range("E13").Value = '  -3.60%  '

$ws.Range("E14").Value = '  +0.40%  '

$ws.Range("D15").Value = '16.36'
$ws.Range("E15").Value = '  +0.19%  '

$ws.Range("D16").Value = '2.728.16'
$ws.Range("E16").Value = '  +0.70%  '

$ws.Range("D17").Value = '2.367.08'
$ws.Range("E17").Value = '  +0.50%  '

$ws.Range("D18").Value = '42.597.99'
$ws.Range("E18").Value = '  -1.14%  '

$ws.Range("D19").Value = '7.89'
$ws.Range("E19").Value = '  +7.94%  '

$ws.Range("E20").Value = '  -0.92%  '

$ws.Range("D21").Value = "'3.80"
$ws.Range("E21").Value = '  +11.77%  '

$ws.Range("D22").Value = '75.76'
$ws.Range("E22").Value = '  +0.73%  '

$ws.Range("D23").Value = '270.69'
$ws.Range("E23").Value = '  +7.29%  '

$ws.Range("D24").Value = '2.32'
$ws.Range("E24").Value = '  -7.48%  '

$ws.Range("D25").Value = '10.01'
$ws.Range("E25").Value = '  +11.94%  '

$ws.Range("D27").Value = '11.52'
$ws.Range("E27").Value = '  -3.93%  '

$ws.Range("D28").Value = '23.38'
$ws.Range("E28").Value = '  +3.80%  '

$ws.Range("E29").Value = '  -1.17%  '

$ws.Range("D30").Value = '175.24'
$ws.Range("E30").Value = '  +1.08%  '

$ws.Range("D31").Value = "'3.10"
$ws.Range("E31").Value = '  -1.89%  '

$ws.Range("D32").Value = '0.0908'
$ws.Range("E32").Value = '  -1.22%  '

$ws.Range("D33").Value = '35.46'
$ws.Range("E33").Value = '  -8.52%  '

$ws.Range("D34").Value = '6.09'
$ws.Range("E34").Value = '  +2.60%  '

$ws.Range("E35").Value = '  +1.53%  '

$ws.Range("E36").Value = '  -8.01%  '

$ws.Range("E37").Value = '  -4.49%  '

$ws.Range("E38").Value = '  +8.31%  '

$ws.Range("B39").Value = 'Kaspa'
$ws.Range("C39").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D39").Value = '0.106'
$ws.Range("E39").Value = '  +3.24%  '

$ws.Range("B40").Value = 'NEARProtocol'
$ws.Range("C40").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D40").Value = '3.84'
$ws.Range("E40").Value = '  -5.75%  '

$ws.Range("E41").Value = '  +2.80%  '

$ws.Range("D42").Value = '0.235'
$ws.Range("E42").Value = '  +1.98%  '

$ws.Range("D43").Value = '70.23'
$ws.Range("E43").Value = '  -2.66%  '

$ws.Range("D45").Value = '117.88'
$ws.Range("E45").Value = '  +7.19%  '

$ws.Range("D46").Value = "'89.20"
$ws.Range("E46").Value = '  +27.99%  '

$ws.Range("D47").Value = '12.03'
$ws.Range("E47").Value = '  -3.10%  '

$ws.Range("E48").Value = '  -2.08%  '

$ws.Range("D49").Value = '9.12'
$ws.Range("E49").Value = '  -2.46%  '

$ws.Range("D50").Value = '1.582.09'
$ws.Range("E50").Value = '  +5.81%  '

$ws.Range("E51").Value = '  -1.40%  '
